$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Price" column header ---
$ws.Range("N1").Value = "Price"

# --- Rows 2-6: blank price cells, just given a plain black font (no data yet) ---
$ws.Range("N2:N6").Font.Color = 0

# --- Row 7: first price entered (kept its own font "Aptos Narrow") ---
$ws.Range("N7").Font.Color = 0
$ws.Range("N7").Font.Name = "Aptos Narrow"
$ws.Range("N7").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("N7").Value = 27.02

# --- Rows 8-50: remaining prices (Calibri font, black, currency format) ---
$ws.Range("N8:N50").Font.Color = 0
$ws.Range("N8:N50").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

$prices = @{
    8  = 30.26
    9  = 31.44
    10 = 39.61
    11 = 39.01
    12 = 34.17
    13 = 37.86
    14 = 44.24
    15 = 34.43
    16 = 45.2
    17 = 30.83
    18 = 27.63
    19 = 18.78
    20 = 20.46
    21 = 14.77
    22 = 8.9700000000000006
    23 = 16.239999999999998
    24 = 11.92
    25 = 25.5
    26 = 23.96
    27 = 22.92
    28 = 30.51
    29 = 33.090000000000003
    30 = 37.58
    31 = 45.64
    32 = 39.47
    33 = 34.799999999999997
    34 = 45.19
    35 = 47.91
    36 = 50.25
    37 = 54.91
    38 = 57.45
    39 = 51.73
    40 = 52.12
    41 = 38.6
    42 = 49.87
    43 = 45.35
    44 = 67.34
    45 = 67.83
    46 = 74.209999999999994
    47 = 99.38
    48 = 86.96
    49 = 83.09
    50 = 88.45
}

foreach ($r in $prices.Keys) {
    $ws.Cells.Item($r, 14).Value = $prices[$r]
}

# --- View/selection bookkeeping to mirror the saved workbook state ---
$ws.Range("N2:N50").Select()
